# Auto-applies the numeric value updates described in the commit diff.
# Each worksheet is addressed by name; cells are updated via .Value,
# new cells are created by assigning .Value, and cells that were removed
# entirely in the diff are cleared with .ClearContents() so the <c> element
# disappears from the saved XML (matching the diff exactly).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 261.5
$ws.Range("I11").Value = 261.5
$ws.Range("K11").Value = 261.5
$ws.Range("M11").Value = -121.5
$ws.Range("H29").Value = 1000
$ws.Range("I29").Value = 1000
$ws.Range("J29").Value = 1000
$ws.Range("K29").Value = 3000
$ws.Range("L29").Value = 3000
$ws.Range("M29").Value = -2719
$ws.Range("N29").Value = -3562
$ws.Range("H64").Value = 6605.6665
$ws.Range("J64").Value = 8387.625
$ws.Range("L64").Value = 8387.625
$ws.Range("N64").Value = -8883.625
$ws.Range("H67").Value = 6605.6665
$ws.Range("J67").Value = 8387.625
$ws.Range("L67").Value = 8387.625
$ws.Range("N67").Value = -10103.625
$ws.Range("H92").Value = 20834066
$ws.Range("I92").Value = 23810256
$ws.Range("K92").Value = 23810256
$ws.Range("M92").Value = -23809008
$ws.Range("H101").Value = 828.8333
$ws.Range("I101").Value = 874
$ws.Range("J101").Value = 603
$ws.Range("K101").Value = 2622
$ws.Range("L101").Value = 1809
$ws.Range("M101").Value = -1000
$ws.Range("N101").Value = -5053
$ws.Range("H106").Value = 16673660
$ws.Range("I106").Value = 16673660
$ws.Range("K106").Value = 16673660
$ws.Range("M106").Value = -16673029
$ws.Range("H116").Value = 14917.059
$ws.Range("I116").Value = 17701.23
$ws.Range("K116").Value = 17701.23
$ws.Range("M116").Value = -14259.23
$ws.Range("H120").Value = 120994.5
$ws.Range("J120").Value = 120994.5
$ws.Range("L120").Value = 120994.5
$ws.Range("N120").Value = -130670.5
$ws.Range("H132").Value = 1409.884
$ws.Range("I132").Value = 1410.4
$ws.Range("K132").Value = 4231.200000000001
$ws.Range("M132").Value = -1701.200000000001
$ws.Range("H137").Value = 1635894.8
$ws.Range("I137").Value = 1666.7894
$ws.Range("K137").Value = 5000.3682
$ws.Range("M137").Value = -2450.3682
$ws.Range("H138").Value = 2744.97
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 2744.97
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 8234.91
$ws.Range("N138").Value = -18514.91
$ws.Range("H141").Value = 3024.7144
$ws.Range("I141").Value = 1862.1666
$ws.Range("K141").Value = 5586.4998
$ws.Range("M141").Value = -406.4997999999996
$ws.Range("M138").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 32018486
$ws.Range("I32").Value = 36380450
$ws.Range("K32").Value = 36380450
$ws.Range("M32").Value = -36380163
$ws.Range("H61").Value = 3365.7932
$ws.Range("I61").Value = 3147.5293
$ws.Range("K61").Value = 3147.5293
$ws.Range("M61").Value = -2935.5293
$ws.Range("H74").Value = 2878.8928
$ws.Range("I74").Value = 2576.36
$ws.Range("J74").Value = 5400
$ws.Range("K74").Value = 2576.36
$ws.Range("L74").Value = 5400
$ws.Range("M74").Value = -1702.36
$ws.Range("N74").Value = -7148
$ws.Range("H77").Value = 2878.8928
$ws.Range("I77").Value = 2576.36
$ws.Range("J77").Value = 5400
$ws.Range("K77").Value = 12881.8
$ws.Range("L77").Value = 27000
$ws.Range("M77").Value = -8513.800000000001
$ws.Range("N77").Value = -35736
$ws.Range("H113").Value = 156000
$ws.Range("J113").Value = 156000
$ws.Range("L113").Value = 156000
$ws.Range("N113").Value = -164678
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("H136").Value = 3365.7932
$ws.Range("I136").Value = 3147.5293
$ws.Range("K136").Value = 9442.5879
$ws.Range("M136").Value = -6892.5879
$ws.Range("H139").Value = 76398.60000000001
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 76398.60000000001
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 76398.60000000001
$ws.Range("N139").Value = -86678.60000000001
$ws.Range("N133").ClearContents()
$ws.Range("M139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2420.9375
$ws.Range("I105").Value = 2057.6956
$ws.Range("K105").Value = 2057.6956
$ws.Range("M105").Value = -310.6956
$ws.Range("H107").Value = 1238.7059
$ws.Range("J107").Value = 1233.0834
$ws.Range("L107").Value = 1233.0834
$ws.Range("N107").Value = -5073.0834

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2811.8333
$ws.Range("I58").Value = 2527.6155
$ws.Range("J58").Value = 4043.4443
$ws.Range("K58").Value = 2527.6155
$ws.Range("L58").Value = 4043.4443
$ws.Range("M58").Value = -2324.6155
$ws.Range("N58").Value = -4449.4443
$ws.Range("H125").Value = 99933
$ws.Range("J125").Value = 99933
$ws.Range("L125").Value = 99933
$ws.Range("N125").Value = -104853
$ws.Range("H132").Value = 4074.6191
$ws.Range("I132").Value = 3987.1765
$ws.Range("K132").Value = 11961.5295
$ws.Range("M132").Value = -9431.529500000001
$ws.Range("H134").Value = 3196.5881
$ws.Range("I134").Value = 3162.8
$ws.Range("J134").Value = 3450
$ws.Range("K134").Value = 9488.400000000001
$ws.Range("L134").Value = 10350
$ws.Range("M134").Value = -6953.400000000001
$ws.Range("N134").Value = -15420
$ws.Range("H136").Value = 2811.8333
$ws.Range("I136").Value = 2527.6155
$ws.Range("J136").Value = 4043.4443
$ws.Range("K136").Value = 7582.8465
$ws.Range("L136").Value = 12130.3329
$ws.Range("M136").Value = -5032.8465
$ws.Range("N136").Value = -17230.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 1134.6666
$ws.Range("I51").Value = 202
$ws.Range("J51").Value = 3000
$ws.Range("K51").Value = 606
$ws.Range("L51").Value = 9000
$ws.Range("M51").Value = -146
$ws.Range("N51").Value = -9920
$ws.Range("H52").Value = 1958.6666
$ws.Range("J52").Value = 1958.6666
$ws.Range("L52").Value = 5875.9998
$ws.Range("N52").Value = -6407.9998
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("H134").Value = 4268.5835
$ws.Range("I134").Value = 4268.5835
$ws.Range("K134").Value = 12805.7505
$ws.Range("M134").Value = -7735.750499999998
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 69966.664
$ws.Range("J64").Value = 69966.664
$ws.Range("L64").Value = 69966.664
$ws.Range("N64").Value = -70462.664
$ws.Range("H67").Value = 69966.664
$ws.Range("J67").Value = 69966.664
$ws.Range("L67").Value = 69966.664
$ws.Range("N67").Value = -71682.664
$ws.Range("H97").Value = 944.8421
$ws.Range("I97").Value = 836.875
$ws.Range("K97").Value = 836.875
$ws.Range("M97").Value = -340.875
$ws.Range("H131").Value = 119000
$ws.Range("J131").Value = 119000
$ws.Range("L131").Value = 119000
$ws.Range("N131").Value = -129080
$ws.Range("H132").Value = 2461.652
$ws.Range("I132").Value = 1830.95
$ws.Range("K132").Value = 5492.85
$ws.Range("M132").Value = -2962.85

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3423.5454
$ws.Range("I46").Value = 411.75
$ws.Range("J46").Value = 4092.8333
$ws.Range("K46").Value = 411.75
$ws.Range("L46").Value = 4092.8333
$ws.Range("M46").Value = -223.75
$ws.Range("N46").Value = -4468.8333
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("H139").Value = 74998
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 22440.4
$ws.Range("I49").Value = 13556
$ws.Range("J49").Value = 28363.334
$ws.Range("K49").Value = 13556
$ws.Range("L49").Value = 28363.334
$ws.Range("M49").Value = -13326
$ws.Range("N49").Value = -28823.334
$ws.Range("H64").Value = 87989.60000000001
$ws.Range("J64").Value = 99982.664
$ws.Range("L64").Value = 99982.664
$ws.Range("N64").Value = -100478.664
$ws.Range("H67").Value = 87989.60000000001
$ws.Range("J67").Value = 99982.664
$ws.Range("L67").Value = 99982.664
$ws.Range("N67").Value = -101698.664
$ws.Range("H74").Value = 63311.5
$ws.Range("J74").Value = 63311.5
$ws.Range("L74").Value = 63311.5
$ws.Range("N74").Value = -65183.5
$ws.Range("H77").Value = 63311.5
$ws.Range("J77").Value = 63311.5
$ws.Range("L77").Value = 189934.5
$ws.Range("N77").Value = -199294.5
$ws.Range("H100").Value = 1165.8125
$ws.Range("I100").Value = 740.1111
$ws.Range("J100").Value = 1713.1428
$ws.Range("K100").Value = 1480.2222
$ws.Range("L100").Value = 3426.2856
$ws.Range("M100").Value = -939.2221999999999
$ws.Range("N100").Value = -4508.2856
$ws.Range("H132").Value = 2637.5186
$ws.Range("I132").Value = 2513.739
$ws.Range("K132").Value = 7541.217000000001
$ws.Range("M132").Value = -5011.217000000001
$ws.Range("H136").Value = 3186.9285
$ws.Range("I136").Value = 2452.25
$ws.Range("J136").Value = 4166.5
$ws.Range("K136").Value = 7356.75
$ws.Range("L136").Value = 12499.5
$ws.Range("M136").Value = -4806.75
$ws.Range("N136").Value = -17599.5
